# "create dc control script"
# Append two new time-sheet entries (rows 42 and 43) below the existing
# logged data (which currently ends at row 41), then move the selection
# down to where the user would naturally continue entering data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 42 - a complete entry (Date, Time In, Time out, + the computed
# Delta Time / Number of minutes / Number of hours / Money formulas that
# every prior row carries).
#
# Grab the number formats from an existing fully-populated row (row 39
# has the same A:G style pattern we need: date / h:mm / h:mm AM-PM /
# h:mm / 0.00 / 0.00 / currency) and paste *just the formatting* onto
# the new row before filling in the real values.
# ---------------------------------------------------------------------
$ws.Range("A39:G39").Copy() | Out-Null
$ws.Range("A42:G42").PasteSpecial(-4122) | Out-Null   # -4122 = xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A42").Value = 45607
$ws.Range("B42").Value = 0.45694444444444443
$ws.Range("C42").Value = 0.6118055555555556
$ws.Range("D42").Formula = "=C42-B42"
$ws.Range("E42").Formula = "=D42*1440"
$ws.Range("F42").Formula = "=E42/60"
$ws.Range("G42").Formula = "=F42*22.5"

# ---------------------------------------------------------------------
# Row 43 - only a clock-in has been recorded so far (Date + Time In);
# Time out / the downstream formulas haven't been entered yet. Copy the
# date/time formatting (style pattern from row 41: date / h:mm AM-PM)
# for just the two populated columns.
# ---------------------------------------------------------------------
$ws.Range("A41:B41").Copy() | Out-Null
$ws.Range("A43:B43").PasteSpecial(-4122) | Out-Null   # -4122 = xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A43").Value = 45607
$ws.Range("B43").Value = 0.62986111111111109

# ---------------------------------------------------------------------
# Leave the cursor where the user ended up after adding these rows.
# ---------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A44").Select() | Out-Null
